# Scheduled runner refresh: re-pull current market prices and recompute
# the Leve profit columns (H..N) on each crafting-class sheet.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 29
$ws.Range("H29").Value = 200.33333
$ws.Range("I29").Value = 200.33333
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 600.99999
$ws.Range("L29").Value = 0
$ws.Range("M29").ClearContents()
$ws.Range("N29").Value = -319.99999
# Row 98
$ws.Range("H98").Value = 1326.1333
$ws.Range("I98").Value = 1475.12
$ws.Range("J98").Value = 581.2
$ws.Range("K98").Value = 1475.12
$ws.Range("L98").Value = 581.2
$ws.Range("M98").Value = 22.88000000000011
$ws.Range("N98").Value = -3577.2
# Row 122
$ws.Range("H122").Value = 1326.1333
$ws.Range("I122").Value = 1475.12
$ws.Range("J122").Value = 581.2
$ws.Range("K122").Value = 4425.36
$ws.Range("L122").Value = 1743.6
$ws.Range("M122").Value = -1975.36
$ws.Range("N122").Value = -6643.6
# Row 129
$ws.Range("H129").Value = 1425455.9
$ws.Range("I129").Value = 477.8
$ws.Range("J129").Value = 1764736.4
$ws.Range("K129").Value = 1433.4
$ws.Range("L129").Value = 5294209.199999999
$ws.Range("M129").Value = 3566.6
$ws.Range("N129").Value = -5304209.199999999
# Row 137
$ws.Range("H137").Value = 992.8889
$ws.Range("I137").Value = 864.8
$ws.Range("J137").Value = 1633.3334
$ws.Range("K137").Value = 2594.4
$ws.Range("L137").Value = 4900.0002
$ws.Range("M137").Value = -44.39999999999964
$ws.Range("N137").Value = -10000.0002

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 9733.879000000001
$ws.Range("I32").Value = 5347.5933
$ws.Range("J32").Value = 46704
$ws.Range("K32").Value = 5347.5933
$ws.Range("L32").Value = 46704
$ws.Range("M32").Value = -5060.5933
$ws.Range("N32").Value = -47278
# Row 82
$ws.Range("H82").Value = 24000
$ws.Range("J82").Value = 24000
$ws.Range("L82").Value = 24000
$ws.Range("N82").Value = -24722
# Row 85
$ws.Range("H85").Value = 24000
$ws.Range("J85").Value = 24000
$ws.Range("L85").Value = 24000
$ws.Range("N85").Value = -26496
# Row 109
$ws.Range("H109").Value = 10000
$ws.Range("J109").Value = 10000
$ws.Range("L109").Value = 10000
$ws.Range("N109").Value = -12774

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 47651516
$ws.Range("I20").Value = 38934.707
$ws.Range("K20").Value = 38934.707
$ws.Range("M20").Value = -38687.707
# Row 81
$ws.Range("H81").Value = 21225
$ws.Range("J81").Value = 21225
$ws.Range("L81").Value = 21225
$ws.Range("N81").Value = -23347
# Row 84
$ws.Range("H84").Value = 21225
$ws.Range("J84").Value = 21225
$ws.Range("L84").Value = 63675
$ws.Range("N84").Value = -74283
# Row 135
$ws.Range("H135").Value = 50259.668
$ws.Range("J135").Value = 50259.668
$ws.Range("L135").Value = 50259.668
$ws.Range("N135").Value = -60399.668
# Row 137
$ws.Range("H137").Value = 46937.5
$ws.Range("J137").Value = 46937.5
$ws.Range("L137").Value = 46937.5
$ws.Range("N137").Value = -57137.5
# Row 138
$ws.Range("H138").Value = 52145
$ws.Range("J138").Value = 52145
$ws.Range("L138").Value = 52145
$ws.Range("N138").Value = -62425

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 344.4
$ws.Range("I22").Value = 342.36365
$ws.Range("J22").Value = 350
$ws.Range("K22").Value = 342.36365
$ws.Range("L22").Value = 350
$ws.Range("M22").Value = 7.636349999999993
$ws.Range("N22").Value = -1050
# Row 31
$ws.Range("H31").Value = 1929.35
$ws.Range("I31").Value = 1873
$ws.Range("K31").Value = 1873
$ws.Range("M31").Value = -1578
# Row 34
$ws.Range("H34").Value = 1929.35
$ws.Range("I34").Value = 1873
$ws.Range("K34").Value = 1873
$ws.Range("M34").Value = -1671
# Row 62
$ws.Range("H62").Value = 8970.714
$ws.Range("I62").Value = 12085.714
$ws.Range("J62").Value = 5855.7144
$ws.Range("K62").Value = 12085.714
$ws.Range("L62").Value = 5855.7144
$ws.Range("M62").Value = -11461.714
$ws.Range("N62").Value = -7103.7144
# Row 65
$ws.Range("H65").Value = 8970.714
$ws.Range("I65").Value = 12085.714
$ws.Range("J65").Value = 5855.7144
$ws.Range("K65").Value = 60428.57
$ws.Range("L65").Value = 29278.572
$ws.Range("M65").Value = -57308.57
$ws.Range("N65").Value = -35518.572
# Row 99
$ws.Range("H99").Value = 2136.0386
$ws.Range("I99").Value = 1627.9375
$ws.Range("J99").Value = 2949
$ws.Range("K99").Value = 1627.9375
$ws.Range("L99").Value = 2949
$ws.Range("M99").Value = -129.9375
$ws.Range("N99").Value = -5945
# Row 107
$ws.Range("H107").Value = 317.82858
$ws.Range("I107").Value = 204.63637
$ws.Range("J107").Value = 369.70834
$ws.Range("K107").Value = 204.63637
$ws.Range("L107").Value = 369.70834
$ws.Range("M107").Value = 1715.36363
$ws.Range("N107").Value = -4209.70834
# Row 126
$ws.Range("H126").Value = 2136.0386
$ws.Range("I126").Value = 1627.9375
$ws.Range("J126").Value = 2949
$ws.Range("K126").Value = 4883.8125
$ws.Range("L126").Value = 8847
$ws.Range("M126").Value = -2413.8125
$ws.Range("N126").Value = -13787
# Row 132
$ws.Range("H132").Value = 1616.0769
$ws.Range("I132").Value = 753.58826
$ws.Range("J132").Value = 3245.2222
$ws.Range("K132").Value = 2260.76478
$ws.Range("L132").Value = 9735.6666
$ws.Range("M132").Value = 269.23522
$ws.Range("N132").Value = -14795.6666

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 1507.4839
$ws.Range("J5").Value = 2860.8333
$ws.Range("L5").Value = 8582.499899999999
$ws.Range("N5").Value = -8806.499899999999
# Row 113
$ws.Range("H113").Value = 805.7143
$ws.Range("I113").Value = 600
$ws.Range("J113").Value = 810
$ws.Range("K113").Value = 1800
$ws.Range("L113").Value = 2430
$ws.Range("M113").Value = 370
$ws.Range("N113").Value = -6770
# Row 135
$ws.Range("H135").Value = 1507.4839
$ws.Range("J135").Value = 2860.8333
$ws.Range("L135").Value = 25747.4997
$ws.Range("N135").Value = -30817.4997

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 7202
$ws.Range("I70").Value = 5404
$ws.Range("K70").Value = 5404
$ws.Range("M70").Value = -5134
# Row 73
$ws.Range("H73").Value = 7202
$ws.Range("I73").Value = 5404
$ws.Range("K73").Value = 5404
$ws.Range("M73").Value = -4468
# Row 97
$ws.Range("H97").Value = 2148.889
$ws.Range("I97").Value = 2424.2856
$ws.Range("J97").Value = 1185
$ws.Range("K97").Value = 2424.2856
$ws.Range("L97").Value = 1185
$ws.Range("M97").Value = -1928.2856
$ws.Range("N97").Value = -2177
# Row 122
$ws.Range("H122").Value = 696199.5
$ws.Range("I122").Value = 1882256.2
$ws.Range("K122").Value = 5646768.6
$ws.Range("M122").Value = -5644318.6
# Row 132
$ws.Range("H132").Value = 2549.0938
$ws.Range("I132").Value = 2168.3914
$ws.Range("J132").Value = 3522
$ws.Range("K132").Value = 6505.174199999999
$ws.Range("L132").Value = 10566
$ws.Range("M132").Value = -3975.174199999999
$ws.Range("N132").Value = -15626
# Row 136
$ws.Range("H136").Value = 22163.25
$ws.Range("J136").Value = 22163.25
$ws.Range("L136").Value = 66489.75
$ws.Range("N136").Value = -71589.75

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 1036.52
$ws.Range("I22").Value = 863.25
$ws.Range("J22").Value = 1196.4615
$ws.Range("K22").Value = 863.25
$ws.Range("L22").Value = 1196.4615
$ws.Range("M22").Value = -568.25
$ws.Range("N22").Value = -1786.4615
# Row 27
$ws.Range("H27").Value = 1036.52
$ws.Range("I27").Value = 863.25
$ws.Range("J27").Value = 1196.4615
$ws.Range("K27").Value = 863.25
$ws.Range("L27").Value = 1196.4615
$ws.Range("M27").Value = -756.25
$ws.Range("N27").Value = -1410.4615
# Row 68
$ws.Range("H68").Value = 2698.5715
$ws.Range("I68").Value = 4266.6665
$ws.Range("J68").Value = 1522.5
$ws.Range("K68").Value = 4266.6665
$ws.Range("L68").Value = 1522.5
$ws.Range("M68").Value = -3517.6665
$ws.Range("N68").Value = -3020.5
# Row 71
$ws.Range("H71").Value = 2698.5715
$ws.Range("I71").Value = 4266.6665
$ws.Range("J71").Value = 1522.5
$ws.Range("K71").Value = 21333.3325
$ws.Range("L71").Value = 7612.5
$ws.Range("M71").Value = -17589.3325
$ws.Range("N71").Value = -15100.5
# Row 100
$ws.Range("H100").Value = 18520852
$ws.Range("I100").Value = 27780228
$ws.Range("K100").Value = 27780228
$ws.Range("M100").Value = -27779687
# Row 132
$ws.Range("H132").Value = 3057.75
$ws.Range("I132").Value = 2737
$ws.Range("J132").Value = 4162.5557
$ws.Range("K132").Value = 8211
$ws.Range("L132").Value = 12487.6671
$ws.Range("M132").Value = -5681
$ws.Range("N132").Value = -17547.6671
# Row 136
$ws.Range("H136").Value = 4277.7856
$ws.Range("I136").Value = 1295.3478
$ws.Range("J136").Value = 17997
$ws.Range("K136").Value = 3886.0434
$ws.Range("L136").Value = 53991
$ws.Range("M136").Value = -1336.0434
$ws.Range("N136").Value = -59091

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 100
$ws.Range("H100").Value = 738.9167
$ws.Range("I100").Value = 303
$ws.Range("J100").Value = 1050.2858
$ws.Range("K100").Value = 606
$ws.Range("L100").Value = 2100.5716
$ws.Range("M100").Value = -65
$ws.Range("N100").Value = -3182.5716
# Row 122
$ws.Range("H122").Value = 1299.8182
$ws.Range("I122").Value = 1237.25
$ws.Range("K122").Value = 3711.75
$ws.Range("M122").Value = -1261.75
# Row 126
$ws.Range("H126").Value = 930
